$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 0) Workbook-level default font -> SimSun (宋体) 12pt, replacing Calibri
# ------------------------------------------------------------------
$wb.Styles.Item("Normal").Font.Name = "宋体"
$wb.Styles.Item("Normal").Font.Size = 12

# ------------------------------------------------------------------
# 1) New English tips column (C) + Chinese key labels already in A
# ------------------------------------------------------------------
$ws.Range("C1").Value = "Please Select"
$ws.Range("C2").Value = "Intranet"
$ws.Range("C3").Value = "EIP that is to be bound "
$ws.Range("C4").Value = "It only supports binding with the EIP that has the full availability zone attribute"

# ------------------------------------------------------------------
# 2) Column B Chinese text (keep/rewrite); rows 3 & 4 need the trailing
#    "IP" in Arial while the rest of the sentence stays in SimSun, so
#    apply the Characters()-level font BEFORE the whole-cell font pass
#    below (that keeps both runs' rPr explicit, matching a real Excel
#    rich-text shared string).
# ------------------------------------------------------------------
$ws.Range("B1").Value = "请选择"
$ws.Range("B2").Value = "内网"

$ws.Range("B3").Value = "要绑定的公网IP"
$b3 = $ws.Range("B3")
$b3.Characters(1, 6).Font.Name = "宋体"
$b3.Characters(1, 6).Font.Size = 12
$b3.Characters(7, 2).Font.Name = "Arial"
$b3.Characters(7, 2).Font.Size = 12

$ws.Range("B4").Value = "仅支持绑定全可用区属性的公网IP"
$b4 = $ws.Range("B4")
$b4.Characters(1, 14).Font.Name = "宋体"
$b4.Characters(1, 14).Font.Size = 12
$b4.Characters(15, 2).Font.Name = "Arial"
$b4.Characters(15, 2).Font.Size = 12

# ------------------------------------------------------------------
# 3) Column D: merged "current selection" helper cell
# ------------------------------------------------------------------
$ws.Range("D1").Value = "请选择内网 xxx 要绑定的公网IP"
$ws.Range("D1:D3").Merge()
$ws.Range("D1:D3").WrapText = $true
$ws.Range("D1").Font.Name = "宋体"
$ws.Range("D1").Font.Size = 12

# ------------------------------------------------------------------
# 4) Keys (col A) and English tips (col C) share the same Arial style,
#    and the rich-text B3/B4 cells fall back to that same cell-level
#    font too (their runs already carry explicit per-character fonts).
# ------------------------------------------------------------------
$arialRange = $ws.Range("A1:A4,C1:C4,B3:B4")
$arialRange.Font.Name = "Arial"
$arialRange.Font.Size = 12

$ws.Range("B1:B2").Font.Name = "宋体"
$ws.Range("B1:B2").Font.Size = 12

# ------------------------------------------------------------------
# 5) Column widths (approximate best-fit layout of the real workbook)
# ------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 18.5
$ws.Columns("B").ColumnWidth = 33
$ws.Columns("C").ColumnWidth = 91.83
$ws.Columns("D").ColumnWidth = 32.67

# ------------------------------------------------------------------
# 6) Page setup (A4 portrait)
# ------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# 7) Selection / view state left on the merged helper cell
# ------------------------------------------------------------------
$ws.Range("D1:D3").Select() | Out-Null

Write-Host "bindPip sheet rebuilt"
